# Weekly update: insert two new price rows (most recent week) at the top of
# the "Plátano" data block (which starts at row 821), pushing the existing
# rows down by two. This matches the behaviour of the author's weekly
# "Fruta / hortaliza" data refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows right before the current row 821.
$ws.Range("A821:A822").EntireRow.Insert()

# Seed both new rows with the constant columns (A,B,C,E-K,Q,R,T) by copying
# the row that is now just below them (the old row 821, shifted to 823).
$ws.Range("A823:T823").Copy($ws.Range("A821:T821"))
$ws.Range("A823:T823").Copy($ws.Range("A822:T822"))

# New row 821: "Pintón" quality entry for the new reporting date.
$ws.Cells.Item(821, 4).Value = 45021    # Fecha
$ws.Cells.Item(821, 12).Value = "Pintón" # Calidad
$ws.Cells.Item(821, 13).Value = 250     # Volumen
$ws.Cells.Item(821, 14).Value = 21000   # Precio mínimo
$ws.Cells.Item(821, 15).Value = 21000   # Precio máximo
$ws.Cells.Item(821, 16).Value = 21000   # Precio promedio ponderado
$ws.Cells.Item(821, 19).Value = 1050    # Precio $/Kg

# New row 822: "Primera Pintón" quality entry for the same new date.
$ws.Cells.Item(822, 4).Value = 45021
$ws.Cells.Item(822, 12).Value = "Primera Pintón"
$ws.Cells.Item(822, 13).Value = 250
$ws.Cells.Item(822, 14).Value = 22000
$ws.Cells.Item(822, 15).Value = 22000
$ws.Cells.Item(822, 16).Value = 22000
$ws.Cells.Item(822, 19).Value = 1100
